# databases/DefaultCONTRACTS.xlsx - refresh order #3's contract details.
# Only the stone_type, phase, and price fields for row 4 actually change;
# everything else in that row stays the same.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Sandstein"   # stone_type: Weißer Stein -> Sandstein
$ws.Range("F4").Value = "Planung"     # phase: Transport -> Planung
$ws.Range("G4").Value = 7500          # price: 1 -> 7500
